$d = $word.ActiveDocument

# Replace the ID placeholder text (consuming the trailing space run too,
# so the two original runs collapse into the single target run).
[void]$d.Content.Find.Execute("**ID__AFFARS_mp_5315_606_90_topic_1__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_MP5315_606_90__ID**", 2)

# Update the first paragraph's formatting: add a (style-less) paragraph
# border that only carries spacing, and widen the left indent.
$p = $d.Paragraphs(1)
$p.Range.ParagraphFormat.Borders.DistanceFromTop = 5
$p.Range.ParagraphFormat.Borders.DistanceFromLeft = 5
$p.Range.ParagraphFormat.Borders.DistanceFromBottom = 5
$p.Range.ParagraphFormat.Borders.DistanceFromRight = 5
$p.Range.ParagraphFormat.LeftIndent = 11.25
